# Update "想去人数" (F column) counts on both the "展览" and "全部类型"
# worksheets, which hold duplicate copies of the same event data.

$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

# Map of row number -> new value for column F
$updates = @{
    3  = 1748
    8  = 12106
    15 = 13530
    16 = 13584
    21 = 1000
    24 = 2044
    25 = 190
}

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $updates.Keys) {
        $ws.Cells.Item($row, 6).Value = $updates[$row]
    }
}
